$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: the highlighted deadline date "Viernes 22 Enero 2022"
# becomes "Viernes 21 Enero 2022" (day 22 -> 21). The matched span sits
# entirely inside the yellow-highlighted run(s), so a plain text
# replacement keeps the existing <w:highlight w:val="yellow"/> run
# formatting intact.
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.ClearFormatting()
$find1.Find.Replacement.ClearFormatting()
$find1.Find.Execute("Viernes 22 Enero 2022", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "Viernes 21 Enero 2022", 2)

# ---------------------------------------------------------------------
# Change 2: "...que no hayan sido ejecutados, durante el mes de
# febrero." becomes "...que no hayan sido ejecutados, durante el mes de
# marzo 2022" (the trailing period is dropped, and "marzo 2022" gets the
# yellow highlight the old "febrero"/"." text never had). We match
# "febrero." (word + following period) so the replacement also consumes
# that final period, and turn on Replacement.Highlight so the new run
# is painted yellow like the other filled-in placeholders.
# ---------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.ClearFormatting()
$find2.Find.Replacement.ClearFormatting()
$find2.Find.Replacement.Highlight = $true
$find2.Find.Execute("febrero.", $false, $false, $false, $false, $false, `
                     $true, 1, $false, "marzo 2022", 2)
